$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,13

$arr[0,0] = 26.45917141421877
$arr[0,1] = 12.911434130392
$arr[0,2] = 4.391986222073093
$arr[0,3] = 9.749158777423457
$arr[0,4] = 50.5670354151373
$arr[0,5] = 0
$arr[0,6] = 7.344005520526261
$arr[0,7] = 34.96074594669072
$arr[0,8] = 9.43458328522118
$arr[0,9] = 0
$arr[0,10] = 12.84096026796656
$arr[0,11] = 0
$arr[0,12] = 20.87467228318588
$arr[1,0] = 26.13633795757178
$arr[1,1] = 12.53158962585921
$arr[1,2] = 4.373720773786421
$arr[1,3] = 9.757491529661475
$arr[1,4] = 50.52868156315164
$arr[1,5] = 0
$arr[1,6] = 7.344005520526261
$arr[1,7] = 35.01450159061546
$arr[1,8] = 9.452025722433971
$arr[1,9] = 0
$arr[1,10] = 12.84170914579773
$arr[1,11] = 0
$arr[1,12] = 20.94969735457067
$arr[2,0] = 25.94309405417679
$arr[2,1] = 12.2965907506034
$arr[2,2] = 4.362163248105512
$arr[2,3] = 9.762969755326878
$arr[2,4] = 50.51795759874373
$arr[2,5] = 0
$arr[2,6] = 7.344005520526261
$arr[2,7] = 35.0552517298282
$arr[2,8] = 9.463303736571183
$arr[2,9] = 0
$arr[2,10] = 12.84432891059917
$arr[2,11] = 0
$arr[2,12] = 20.99778645518974
$arr[3,0] = 25.86568777793888
$arr[3,1] = 12.20054596235998
$arr[3,2] = 4.357366874866108
$arr[3,3] = 9.765293369496026
$arr[3,4] = 50.51680918569122
$arr[3,5] = 0
$arr[3,6] = 7.344005520526261
$arr[3,7] = 35.07379794977528
$arr[3,8] = 9.468042985702652
$arr[3,9] = 0
$arr[3,10] = 12.84594017774414
$arr[3,11] = 0
$arr[3,12] = 21.01789355760117
$arr[4,0] = 25.85291818287794
$arr[4,1] = 12.1845858274531
$arr[4,2] = 4.356565197219108
$arr[4,3] = 9.765684717443971
$arr[4,4] = 50.5168129293349
$arr[4,5] = 0
$arr[4,6] = 7.344005520526261
$arr[4,7] = 35.07699450221219
$arr[4,8] = 9.468838607121528
$arr[4,9] = 0
$arr[4,10] = 12.84624057737374
$arr[4,11] = 0
$arr[4,12] = 21.02126320082472
$arr[5,0] = 25.94204457405209
$arr[5,1] = 12.29529636965991
$arr[5,2] = 4.362098913725326
$arr[5,3] = 9.76300072289127
$arr[5,4] = 50.51792907283059
$arr[5,5] = 0
$arr[5,6] = 7.344005520526261
$arr[5,7] = 35.05549400472929
$arr[5,8] = 9.463367070641073
$arr[5,9] = 0
$arr[5,10] = 12.84434843875106
$arr[5,11] = 0
$arr[5,12] = 20.99805555794628
$arr[6,0] = 26.34688061612836
$arr[6,1] = 12.78092785368479
$arr[6,2] = 4.38575908859574
$arr[6,3] = 9.751956918173866
$arr[6,4] = 50.55114774592538
$arr[6,5] = 0
$arr[6,6] = 7.344005520526261
$arr[6,7] = 34.97766952846975
$arr[6,8] = 9.440479770867288
$arr[6,9] = 0
$arr[6,10] = 12.84077057065546
$arr[6,11] = 0
$arr[6,12] = 20.90012194028181
$arr[7,0] = 27.17615630022571
$arr[7,1] = 13.71241667646494
$arr[7,2] = 4.429461493143265
$arr[7,3] = 9.733162822080303
$arr[7,4] = 50.71810602218076
$arr[7,5] = 0
$arr[7,6] = 7.344005520526261
$arr[7,7] = 34.88679777846826
$arr[7,8] = 9.400085652749688
$arr[7,9] = 0
$arr[7,10] = 12.85086184844267
$arr[7,11] = 0
$arr[7,12] = 20.72405446489119
$arr[8,0] = 27.80140797251802
$arr[8,1] = 14.37575874740169
$arr[8,2] = 4.459958681125061
$arr[8,3] = 9.721088383441003
$arr[8,4] = 50.90276927986714
$arr[8,5] = 0
$arr[8,6] = 7.344005520526261
$arr[8,7] = 34.85805098449552
$arr[8,8] = 9.373114147156397
$arr[8,9] = 0
$arr[8,10] = 12.86865429925096
$arr[8,11] = 0
$arr[8,12] = 20.60433261798779
$arr[9,0] = 28.0881233586453
$arr[9,1] = 14.6714271889891
$arr[9,2] = 4.4734897488278
$arr[9,3] = 9.715969427790242
$arr[9,4] = 51.0001732855513
$arr[9,5] = 0
$arr[9,6] = 7.344005520526261
$arr[9,7] = 34.8532971450586
$arr[9,8] = 9.361425346720671
$arr[9,9] = 0
$arr[9,10] = 12.87898710214158
$arr[9,11] = 0
$arr[9,12] = 20.55193731565155
$arr[10,0] = 28.1969241904007
$arr[10,1] = 14.78240060999087
$arr[10,2] = 4.478564908464571
$arr[10,3] = 9.71408457712637
$arr[10,4] = 51.03897447646676
$arr[10,5] = 0
$arr[10,6] = 7.344005520526261
$arr[10,7] = 34.85269835654348
$arr[10,8] = 9.357082118524787
$arr[10,9] = 0
$arr[10,10] = 12.88322013830566
$arr[10,11] = 0
$arr[10,12] = 20.53239209602322
$arr[11,0] = 28.17348343299082
$arr[11,1] = 14.75854622074576
$arr[11,2] = 4.477474043527016
$arr[11,3] = 9.714488132677733
$arr[11,4] = 51.03053290636326
$arr[11,5] = 0
$arr[11,6] = 7.344005520526261
$arr[11,7] = 34.85277381948808
$arr[11,8] = 9.358013822874538
$arr[11,9] = 0
$arr[11,10] = 12.88229426582746
$arr[11,11] = 0
$arr[11,12] = 20.5365883784639
$arr[12,0] = 28.09707047140583
$arr[12,1] = 14.68057749099452
$arr[12,2] = 4.473908259931508
$arr[12,3] = 9.715813286939518
$arr[12,4] = 51.00332714412053
$arr[12,5] = 0
$arr[12,6] = 7.344005520526261
$arr[12,7] = 34.85322378128879
$arr[12,8] = 9.361066364411791
$arr[12,9] = 0
$arr[12,10] = 12.87932894939872
$arr[12,11] = 0
$arr[12,12] = 20.55032339979959
$arr[13,0] = 28.05029204877494
$arr[13,1] = 14.63268729129243
$arr[13,2] = 4.471717777719594
$arr[13,3] = 9.716631956063472
$arr[13,4] = 50.98691204843453
$arr[13,5] = 0
$arr[13,6] = 7.344005520526261
$arr[13,7] = 34.85365597461587
$arr[13,8] = 9.362946939188934
$arr[13,9] = 0
$arr[13,10] = 12.87755425560843
$arr[13,11] = 0
$arr[13,12] = 20.55877496865298
$arr[14,0] = 27.782707793034
$arr[14,1] = 14.35630357535201
$arr[14,2] = 4.459067545798236
$arr[14,3] = 9.721430426579122
$arr[14,4] = 50.89667250003298
$arr[14,5] = 0
$arr[14,6] = 7.344005520526261
$arr[14,7] = 34.85852954502899
$arr[14,8] = 9.373889685281391
$arr[14,9] = 0
$arr[14,10] = 12.86802393658912
$arr[14,11] = 0
$arr[14,12] = 20.60779823820613
$arr[15,0] = 27.6190667517188
$arr[15,1] = 14.18510675610438
$arr[15,2] = 4.451219701717309
$arr[15,3] = 9.724469748198802
$arr[15,4] = 50.84473882236504
$arr[15,5] = 0
$arr[15,6] = 7.344005520526261
$arr[15,7] = 34.86365452496313
$arr[15,8] = 9.380751120683719
$arr[15,9] = 0
$arr[15,10] = 12.86274967831302
$arr[15,11] = 0
$arr[15,12] = 20.63840073544576
$arr[16,0] = 27.52516580555331
$arr[16,1] = 14.08607184114661
$arr[16,2] = 4.446673587546248
$arr[16,3] = 9.726253075709987
$arr[16,4] = 50.81612998235783
$arr[16,5] = 0
$arr[16,6] = 7.344005520526261
$arr[16,7] = 34.86738547070256
$arr[16,8] = 9.384752318254378
$arr[16,9] = 0
$arr[16,10] = 12.85992692025499
$arr[16,11] = 0
$arr[16,12] = 20.65619710242613
$arr[17,0] = 27.49341361576404
$arr[17,1] = 14.05244659343968
$arr[17,2] = 4.445128799998725
$arr[17,3] = 9.726862928648236
$arr[17,4] = 50.80666055396026
$arr[17,5] = 0
$arr[17,6] = 7.344005520526261
$arr[17,7] = 34.8687830702373
$arr[17,8] = 9.386116459822762
$arr[17,9] = 0
$arr[17,10] = 12.85900744658431
$arr[17,11] = 0
$arr[17,12] = 20.66225611715946
$arr[18,0] = 27.6364644759
$arr[18,1] = 14.20339049778677
$arr[18,2] = 4.452058451728838
$arr[18,3] = 9.724142566356536
$arr[18,4] = 50.85013670267138
$arr[18,5] = 0
$arr[18,6] = 7.344005520526261
$arr[18,7] = 34.86302786979218
$arr[18,8] = 9.380015053002579
$arr[18,9] = 0
$arr[18,10] = 12.8632893211571
$arr[18,11] = 0
$arr[18,12] = 20.63512291553068
$arr[19,0] = 28.11950941885348
$arr[19,1] = 14.70350651508582
$arr[19,2] = 4.47495693678631
$arr[19,3] = 9.715422603943608
$arr[19,4] = 51.01126622108488
$arr[19,5] = 0
$arr[19,6] = 7.344005520526261
$arr[19,7] = 34.85305897859917
$arr[19,8] = 9.360167507995465
$arr[19,9] = 0
$arr[19,10] = 12.88019125866601
$arr[19,11] = 0
$arr[19,12] = 20.5462810761944
$arr[20,0] = 28.4364903492973
$arr[20,1] = 15.02454666735179
$arr[20,2] = 4.489638407943154
$arr[20,3] = 9.710035877064289
$arr[20,4] = 51.12773734243503
$arr[20,5] = 0
$arr[20,6] = 7.344005520526261
$arr[20,7] = 34.853547956109
$arr[20,8] = 9.347679968060211
$arr[20,9] = 0
$arr[20,10] = 12.89310323051219
$arr[20,11] = 0
$arr[20,12] = 20.48994095713617
$arr[21,0] = 28.26722707845249
$arr[21,1] = 14.85376841875835
$arr[21,2] = 4.48182846796245
$arr[21,3] = 9.71288235370899
$arr[21,4] = 51.06455714276804
$arr[21,5] = 0
$arr[21,6] = 7.344005520526261
$arr[21,7] = 34.85264477363728
$arr[21,8] = 9.354300661518684
$arr[21,9] = 0
$arr[21,10] = 12.8860417944705
$arr[21,11] = 0
$arr[21,12] = 20.51985354667939
$arr[22,0] = 27.62859840125619
$arr[22,1] = 14.19512631433413
$arr[22,2] = 4.451679359669735
$arr[22,3] = 9.724290373182304
$arr[22,4] = 50.84769243135194
$arr[22,5] = 0
$arr[22,6] = 7.344005520526261
$arr[22,7] = 34.86330873732074
$arr[22,8] = 9.380347653286584
$arr[22,9] = 0
$arr[22,10] = 12.86304469600534
$arr[22,11] = 0
$arr[22,12] = 20.6366041868014
$arr[23,0] = 26.94860057769157
$arr[23,1] = 13.46355644556004
$arr[23,2] = 4.417924044896679
$arr[23,3] = 9.737941859395388
$arr[23,4] = 50.66203660410905
$arr[23,5] = 0
$arr[23,6] = 7.344005520526261
$arr[23,7] = 34.90472982267066
$arr[23,8] = 9.410535975687115
$arr[23,9] = 0
$arr[23,10] = 12.84630481507732
$arr[23,11] = 0
$arr[23,12] = 20.76998570992131

$ws.Range("B2:N25").Value = $arr
